# Repull data, push all data, mean calculation
# Updates the dSF column (F) values for the rows whose underlying data was
# re-pulled, reflecting the recalculated deltas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -7
    5  = 1
    11 = -2
    12 = -5
    15 = 2
    16 = 3
    18 = 2
    20 = 0
    23 = -1
    27 = -1
    37 = -1
    38 = 1
    39 = 0
    44 = -2
    51 = 1
    56 = 1
    62 = 0
    65 = -2
    66 = -3
    67 = 1
    68 = 4
    70 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
